# Tighten the hanging-bullet indent on the resume's "No Spacing" bullet
# lists: w:ind w:left="450" (no hanging)  ->  w:ind w:left="360" w:hanging="270"
#
# In Word's object model, left indent / hanging indent are expressed in
# points rather than twips (1 pt = 20 twips), and a hanging indent is a
# negative FirstLineIndent:
#   450 twips = 22.5 pt   (old LeftIndent)
#   360 twips = 18 pt     (new LeftIndent)
#   270 twips = 13.5 pt   (new hanging amount -> FirstLineIndent = -13.5)

$d = $word.ActiveDocument

$oldLeft = 22.5
$newLeft = 18
$newFirstLine = -13.5
$tolerance = 0.01

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ([Math]::Abs($p.LeftIndent - $oldLeft) -lt $tolerance) {
        $p.LeftIndent = $newLeft
        $p.Format.FirstLineIndent = $newFirstLine
    }
}
